$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.019476183077629
$ws.Cells.Item(2, 4).Value = 1.025677550897424
$ws.Cells.Item(2, 5).Value = 1.020606657726835
$ws.Cells.Item(2, 6).Value = 1.03351062202601
$ws.Cells.Item(2, 9).Value = 1.029865279986646
$ws.Cells.Item(2, 10).Value = 1.024678733671311
$ws.Cells.Item(2, 11).Value = 1.028502421198872
$ws.Cells.Item(2, 12).Value = 1.023446444933005
$ws.Cells.Item(2, 13).Value = 1.036312762704638
$ws.Cells.Item(2, 14).Value = 1.026133895005419

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.020446586521331
$ws.Cells.Item(3, 4).Value = 1.026424952079108
$ws.Cells.Item(3, 5).Value = 1.021430123867957
$ws.Cells.Item(3, 6).Value = 1.035716949087132
$ws.Cells.Item(3, 9).Value = 1.030133699082489
$ws.Cells.Item(3, 10).Value = 1.025285799240342
$ws.Cells.Item(3, 11).Value = 1.0290573387006
$ws.Cells.Item(3, 12).Value = 1.024076139128947
$ws.Cells.Item(3, 13).Value = 1.038324360259656
$ws.Cells.Item(3, 14).Value = 1.026741822677188

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.021074449065226
$ws.Cells.Item(4, 4).Value = 1.026908027233316
$ws.Cells.Item(4, 5).Value = 1.021963307740601
$ws.Cells.Item(4, 6).Value = 1.037138209544582
$ws.Cells.Item(4, 9).Value = 1.030305048855544
$ws.Cells.Item(4, 10).Value = 1.025677987872938
$ws.Cells.Item(4, 11).Value = 1.029415204152155
$ws.Cells.Item(4, 12).Value = 1.024483295877081
$ws.Cells.Item(4, 13).Value = 1.039619299230912
$ws.Cells.Item(4, 14).Value = 1.027134568262627

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.021338390393287
$ws.Cells.Item(5, 4).Value = 1.02711098247669
$ws.Cells.Item(5, 5).Value = 1.022187541146269
$ws.Cells.Item(5, 6).Value = 1.037734214019464
$ws.Cells.Item(5, 9).Value = 1.030376527166899
$ws.Cells.Item(5, 10).Value = 1.025842715198604
$ws.Cells.Item(5, 11).Value = 1.029565363504006
$ws.Cells.Item(5, 12).Value = 1.024654393609589
$ws.Cells.Item(5, 13).Value = 1.040162119569517
$ws.Cells.Item(5, 14).Value = 1.02729952951999

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.021382706614649
$ws.Cells.Item(6, 4).Value = 1.02714505197304
$ws.Cells.Item(6, 5).Value = 1.022225195783537
$ws.Cells.Item(6, 6).Value = 1.037834199200496
$ws.Cells.Item(6, 9).Value = 1.030388496075135
$ws.Cells.Item(6, 10).Value = 1.025870364930161
$ws.Cells.Item(6, 11).Value = 1.029590559095556
$ws.Cells.Item(6, 12).Value = 1.024683117517441
$ws.Cells.Item(6, 13).Value = 1.040253170224409
$ws.Cells.Item(6, 14).Value = 1.027327218517338

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.021077975908347
$ws.Cells.Item(7, 4).Value = 1.026910739639233
$ws.Cells.Item(7, 5).Value = 1.021966303632235
$ws.Cells.Item(7, 6).Value = 1.037146179207844
$ws.Cells.Item(7, 9).Value = 1.030306006138628
$ws.Cells.Item(7, 10).Value = 1.025680189551229
$ws.Cells.Item(7, 11).Value = 1.029417211716291
$ws.Cells.Item(7, 12).Value = 1.024485582371591
$ws.Cells.Item(7, 13).Value = 1.03962655855705
$ws.Cells.Item(7, 14).Value = 1.027136773067554

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.019804146940368
$ws.Cells.Item(8, 4).Value = 1.025930251472447
$ws.Cells.Item(8, 5).Value = 1.020884880298615
$ws.Cells.Item(8, 6).Value = 1.034257602436111
$ws.Cells.Item(8, 9).Value = 1.02995647820935
$ws.Cells.Item(8, 10).Value = 1.024884023685397
$ws.Cells.Item(8, 11).Value = 1.028690207795383
$ws.Cells.Item(8, 12).Value = 1.02365931479738
$ws.Cells.Item(8, 13).Value = 1.036993997891865
$ws.Cells.Item(8, 14).Value = 1.026339476554871

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.017559056061813
$ws.Cells.Item(9, 4).Value = 1.024198317979585
$ws.Cells.Item(9, 5).Value = 1.018981923107847
$ws.Cells.Item(9, 6).Value = 1.029117144706772
$ws.Cells.Item(9, 9).Value = 1.02932258993816
$ws.Cells.Item(9, 10).Value = 1.023476280904734
$ws.Cells.Item(9, 11).Value = 1.027399881174223
$ws.Cells.Item(9, 12).Value = 1.022201034104683
$ws.Cells.Item(9, 13).Value = 1.0323023399405
$ws.Cells.Item(9, 14).Value = 1.024929734617989

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016061981714906
$ws.Cells.Item(10, 4).Value = 1.023040837776893
$ws.Cells.Item(10, 5).Value = 1.017715057668187
$ws.Cells.Item(10, 6).Value = 1.025654089861128
$ws.Cells.Item(10, 9).Value = 1.028887783089274
$ws.Cells.Item(10, 10).Value = 1.022534525388574
$ws.Cells.Item(10, 11).Value = 1.026533385995896
$ws.Cells.Item(10, 12).Value = 1.021227284021351
$ws.Cells.Item(10, 13).Value = 1.029137065030688
$ws.Cells.Item(10, 14).Value = 1.023986641700975

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015413637093642
$ws.Cells.Item(11, 4).Value = 1.022538947287413
$ws.Cells.Item(11, 5).Value = 1.017166908190815
$ws.Cells.Item(11, 6).Value = 1.02414549782527
$ws.Cells.Item(11, 9).Value = 1.028696579903233
$ws.Cells.Item(11, 10).Value = 1.022125952521608
$ws.Cells.Item(11, 11).Value = 1.026156681210696
$ws.Cells.Item(11, 12).Value = 1.020805262675869
$ws.Cells.Item(11, 13).Value = 1.027757111595716
$ws.Cells.Item(11, 14).Value = 1.023577488613674

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015172796557848
$ws.Cells.Item(12, 4).Value = 1.022352417470139
$ws.Cells.Item(12, 5).Value = 1.016963362507265
$ws.Cells.Item(12, 6).Value = 1.02358373461059
$ws.Cells.Item(12, 9).Value = 1.028625115795893
$ws.Cells.Item(12, 10).Value = 1.021974071308858
$ws.Cells.Item(12, 11).Value = 1.026016528595134
$ws.Cells.Item(12, 12).Value = 1.020648447259912
$ws.Cells.Item(12, 13).Value = 1.027243087983141
$ws.Cells.Item(12, 14).Value = 1.023425391712181

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015224458435847
$ws.Cells.Item(13, 4).Value = 1.022392433534128
$ws.Cells.Item(13, 5).Value = 1.017007020973919
$ws.Cells.Item(13, 6).Value = 1.023704299000578
$ws.Cells.Item(13, 9).Value = 1.028640465167855
$ws.Cells.Item(13, 10).Value = 1.022006655746916
$ws.Cells.Item(13, 11).Value = 1.026046602131997
$ws.Cells.Item(13, 12).Value = 1.020682087329508
$ws.Cells.Item(13, 13).Value = 1.027353413942931
$ws.Cells.Item(13, 14).Value = 1.023458022423879

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.015393729483106
$ws.Cells.Item(14, 4).Value = 1.022523530829972
$ws.Cells.Item(14, 5).Value = 1.017150081790235
$ws.Cells.Item(14, 6).Value = 1.024099091214166
$ws.Cells.Item(14, 9).Value = 1.028690681705189
$ws.Cells.Item(14, 10).Value = 1.022113400407231
$ws.Cells.Item(14, 11).Value = 1.026145100804343
$ws.Cells.Item(14, 12).Value = 1.020792301449563
$ws.Cells.Item(14, 13).Value = 1.027714652005192
$ws.Cells.Item(14, 14).Value = 1.023564918673854

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015498020741676
$ws.Cells.Item(15, 4).Value = 1.022604290208312
$ws.Cells.Item(15, 5).Value = 1.017238234409383
$ws.Cells.Item(15, 6).Value = 1.024342148260511
$ws.Cells.Item(15, 9).Value = 1.028721563023522
$ws.Cells.Item(15, 10).Value = 1.022179153511761
$ws.Cells.Item(15, 11).Value = 1.026205758876155
$ws.Cells.Item(15, 12).Value = 1.020860200314846
$ws.Cells.Item(15, 13).Value = 1.027937029686626
$ws.Cells.Item(15, 14).Value = 1.023630765155336

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016105007969346
$ws.Cells.Item(16, 4).Value = 1.023074131883803
$ws.Cells.Item(16, 5).Value = 1.017751445182461
$ws.Cells.Item(16, 6).Value = 1.025754015488661
$ws.Cells.Item(16, 9).Value = 1.028900410677338
$ws.Cells.Item(16, 10).Value = 1.022561624377748
$ws.Cells.Item(16, 11).Value = 1.02655835479611
$ws.Cells.Item(16, 12).Value = 1.021255284104481
$ws.Cells.Item(16, 13).Value = 1.029228447218378
$ws.Cells.Item(16, 14).Value = 1.024013779173822

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.016485726982882
$ws.Cells.Item(17, 4).Value = 1.023368664617649
$ws.Cells.Item(17, 5).Value = 1.018073478669456
$ws.Cells.Item(17, 6).Value = 1.026637185166668
$ws.Cells.Item(17, 9).Value = 1.029011811043871
$ws.Cells.Item(17, 10).Value = 1.02280132705504
$ws.Cells.Item(17, 11).Value = 1.026779124604173
$ws.Cells.Item(17, 12).Value = 1.021503007247487
$ws.Cells.Item(17, 13).Value = 1.030035983128468
$ws.Cells.Item(17, 14).Value = 1.024253822256411

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.016707784229648
$ws.Cells.Item(18, 4).Value = 1.023540393702046
$ws.Cells.Item(18, 5).Value = 1.018261355099393
$ws.Cells.Item(18, 6).Value = 1.02715145125065
$ws.Cells.Item(18, 9).Value = 1.02907650657708
$ws.Cells.Item(18, 10).Value = 1.022941065747235
$ws.Cells.Item(18, 11).Value = 1.026907750620852
$ws.Cells.Item(18, 12).Value = 1.021647463237453
$ws.Cells.Item(18, 13).Value = 1.030506103972261
$ws.Cells.Item(18, 14).Value = 1.024393759393576

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.016783498400758
$ws.Cells.Item(19, 4).Value = 1.023598937570322
$ws.Cells.Item(19, 5).Value = 1.018325422859546
$ws.Cells.Item(19, 6).Value = 1.027326656099642
$ws.Cells.Item(19, 9).Value = 1.029098518259598
$ws.Cells.Item(19, 10).Value = 1.022988700195174
$ws.Cells.Item(19, 11).Value = 1.026951584166106
$ws.Cells.Item(19, 12).Value = 1.021696712785632
$ws.Cells.Item(19, 13).Value = 1.030666251371573
$ws.Cells.Item(19, 14).Value = 1.024441461487895

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.016444880449385
$ws.Cells.Item(20, 4).Value = 1.023337070967965
$ws.Cells.Item(20, 5).Value = 1.018038923415953
$ws.Cells.Item(20, 6).Value = 1.026542519905386
$ws.Cells.Item(20, 9).Value = 1.028999888069771
$ws.Cells.Item(20, 10).Value = 1.022775617075642
$ws.Cells.Item(20, 11).Value = 1.026755453129435
$ws.Cells.Item(20, 12).Value = 1.021476432702313
$ws.Cells.Item(20, 13).Value = 1.029949435669071
$ws.Cells.Item(20, 14).Value = 1.024228075765893

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015343883843439
$ws.Cells.Item(21, 4).Value = 1.022484928860829
$ws.Cells.Item(21, 5).Value = 1.017107952252913
$ws.Cells.Item(21, 6).Value = 1.023982873878937
$ws.Cells.Item(21, 9).Value = 1.028675906430006
$ws.Cells.Item(21, 10).Value = 1.022081970059246
$ws.Cells.Item(21, 11).Value = 1.026116101696484
$ws.Cells.Item(21, 12).Value = 1.020759847741751
$ws.Cells.Item(21, 13).Value = 1.027608316681085
$ws.Cells.Item(21, 14).Value = 1.023533443691171

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014651547624334
$ws.Cells.Item(22, 4).Value = 1.021948543239114
$ws.Cells.Item(22, 5).Value = 1.016522968973869
$ws.Cells.Item(22, 6).Value = 1.02236537235518
$ws.Cells.Item(22, 9).Value = 1.028469643431642
$ws.Cells.Item(22, 10).Value = 1.021645157346518
$ws.Cells.Item(22, 11).Value = 1.025712798063786
$ws.Cells.Item(22, 12).Value = 1.020308967018603
$ws.Cells.Item(22, 13).Value = 1.026127967340579
$ws.Cells.Item(22, 14).Value = 1.023096010654287

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015018577664565
$ws.Cells.Item(23, 4).Value = 1.022232949608591
$ws.Cells.Item(23, 5).Value = 1.016833046148152
$ws.Cells.Item(23, 6).Value = 1.023223627595889
$ws.Cells.Item(23, 9).Value = 1.028579231195812
$ws.Cells.Item(23, 10).Value = 1.021876785675314
$ws.Cells.Item(23, 11).Value = 1.025926722337779
$ws.Cells.Item(23, 12).Value = 1.020548019425556
$ws.Cells.Item(23, 13).Value = 1.026913537808795
$ws.Cells.Item(23, 14).Value = 1.023327967921879

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016463337270131
$ws.Cells.Item(24, 4).Value = 1.023351346985338
$ws.Cells.Item(24, 5).Value = 1.018054537325214
$ws.Cells.Item(24, 6).Value = 1.026585297754383
$ws.Cells.Item(24, 9).Value = 1.029005276421248
$ws.Cells.Item(24, 10).Value = 1.022787234543259
$ws.Cells.Item(24, 11).Value = 1.02676614969911
$ws.Cells.Item(24, 12).Value = 1.021488440709203
$ws.Cells.Item(24, 13).Value = 1.029988545523612
$ws.Cells.Item(24, 14).Value = 1.024239709731646

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.018139523248009
$ws.Cells.Item(25, 4).Value = 1.02464656477325
$ws.Cells.Item(25, 5).Value = 1.019473569989327
$ws.Cells.Item(25, 6).Value = 1.030452273543295
$ws.Cells.Item(25, 9).Value = 1.029488608508511
$ws.Cells.Item(25, 10).Value = 1.023840788115297
$ws.Cells.Item(25, 11).Value = 1.027734563689688
$ws.Cells.Item(25, 12).Value = 1.022578308370808
$ws.Cells.Item(25, 13).Value = 1.0335217046009
$ws.Cells.Item(25, 14).Value = 1.0252947594706

Write-Output "Applied vm_pu updates for rows 2-25"